# Update Receptor expression values (Target cluster = ECs, rows 2 and 5)
# and the downstream derived-specificity / edge-weight columns that are
# recomputed from them, reflecting the new TPM-based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (FAPs -> Angptl3/Itgav -> ECs) ---
$ws.Range("M2").Value = 3.759736666666667
$ws.Range("N2").Value = 11.27921
$ws.Range("O2").Value = 0.0683751702595819
$ws.Range("P2").Value = 0.06837517025958188
$ws.Range("Q2").Value = 8.678224173999999
$ws.Range("R2").Value = 78.104017566
$ws.Range("S2").Value = 0.04538941146465603
$ws.Range("T2").Value = 0.04538941146465601

# --- Row 3 (FAPs -> Angptl3/Itgav -> FAPs) : O/P/S/T renormalized ---
$ws.Range("O3").Value = 0.6514180024294648
$ws.Range("P3").Value = 0.6514180024294647
$ws.Range("S3").Value = 0.4324300712598485
$ws.Range("T3").Value = 0.4324300712598484

# --- Row 4 (FAPs -> Angptl3/Itgav -> MuSCs) : O/P/S/T renormalized ---
$ws.Range("O4").Value = 0.2802068273109533
$ws.Range("P4").Value = 0.2802068273109533
$ws.Range("S4").Value = 0.1860093793074007
$ws.Range("T4").Value = 0.1860093793074007

# --- Row 5 (MuSCs -> Angptl3/Itgav -> ECs) ---
$ws.Range("M5").Value = 3.759736666666667
$ws.Range("N5").Value = 11.27921
$ws.Range("O5").Value = 0.0683751702595819
$ws.Range("P5").Value = 0.06837517025958188
$ws.Range("Q5").Value = 4.394759949403333
$ws.Range("R5").Value = 39.55283954463
$ws.Range("S5").Value = 0.02298575879492587
$ws.Range("T5").Value = 0.02298575879492586

# --- Row 6 (MuSCs -> Angptl3/Itgav -> FAPs) : O/P/S/T renormalized ---
$ws.Range("O6").Value = 0.6514180024294648
$ws.Range("P6").Value = 0.6514180024294647
$ws.Range("S6").Value = 0.2189879311696163
$ws.Range("T6").Value = 0.2189879311696162

# --- Row 7 (MuSCs -> Angptl3/Itgav -> MuSCs) : O/P/S/T renormalized ---
$ws.Range("O7").Value = 0.2802068273109533
$ws.Range("P7").Value = 0.2802068273109533
$ws.Range("S7").Value = 0.09419744800355258
$ws.Range("T7").Value = 0.09419744800355255
